$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns for every coin row
# with the latest scraped figures. Several "Price" values look like plain
# numbers (e.g. "1.004", "320.59") but must stay literal TEXT, exactly as
# they were authored (matching the original inlineStr cells) - a leading
# apostrophe forces Excel to store them as text instead of silently
# re-parsing them as floats (which would also mangle trailing zeros such
# as "3.600" -> 3.6 or "0.00001020" -> 1.02E-05).
$ws.Range("D2").Value = "27.387.08"
$ws.Range("D3").Value = "1.849.18"
$ws.Range("E3").Value = "  -5.52%  "
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").Value = "'320.59"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "'0.4456"
$ws.Range("E7").Value = "  -6.47%  "
$ws.Range("D8").Value = "'0.3831"
$ws.Range("E8").Value = "  -5.74%  "
$ws.Range("D9").Value = "'48.37"
$ws.Range("E9").Value = "  -9.15%  "
$ws.Range("D10").Value = "'0.07803"
$ws.Range("E10").Value = "  -7.60%  "
$ws.Range("D11").Value = "'1.014"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").Value = "'21.47"
$ws.Range("E12").Value = "  -3.11%  "
$ws.Range("D13").Value = "1.844.06"
$ws.Range("E13").Value = "  -6.63%  "
$ws.Range("D14").Value = "'5.827"
$ws.Range("E14").Value = "  -5.95%  "
$ws.Range("D15").Value = "'7.081"
$ws.Range("E15").Value = "  -7.25%  "
$ws.Range("D16").Value = "'1.006"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'85.43"
$ws.Range("E17").Value = "  -4.32%  "
$ws.Range("D18").Value = "'0.00001020"
$ws.Range("E18").Value = "  -5.09%  "
$ws.Range("D19").Value = "'0.06495"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "'16.91"
$ws.Range("E20").Value = "  -9.75%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'5.459"
$ws.Range("E22").Value = "  -6.21%  "
$ws.Range("D23").Value = "27.381.43"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").Value = "'10.73"
$ws.Range("E24").Value = "  -7.60%  "
$ws.Range("D25").Value = "'2.273"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "2.070.99"
$ws.Range("E26").Value = "  -6.12%  "
$ws.Range("D27").Value = "'151.37"
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("D28").Value = "'19.29"
$ws.Range("E28").Value = "  -4.73%  "

# Rows 29/30 swapped rank order: LidoDAOToken now outranks
# InternetComputer(DFINITY), so the two rows' Coin/Link/Price/Volume
# contents trade places (the rank numbers in column A stay put).
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.023"
$ws.Range("E29").Value = "  -6.61%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'5.450"
$ws.Range("E30").Value = "  -9.24%  "
$ws.Range("D31").Value = "'119.38"
$ws.Range("E31").Value = "  -3.67%  "
$ws.Range("D32").Value = "'1.482"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "'0.09307"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "'0.9225"
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("D35").Value = "'3.600"
$ws.Range("E35").Value = "  -1.65%  "
$ws.Range("D36").Value = "'5.192"
$ws.Range("E36").Value = "  -7.36%  "
$ws.Range("D37").Value = "'0.02210"
$ws.Range("E37").Value = "  -5.53%  "

# Rows 38/39 likewise swapped: Hedera now outranks TrustWalletToken.
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05956"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.208"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("D40").Value = "'8.271"
$ws.Range("E40").Value = "  -6.01%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("D42").Value = "'0.5874"
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("D43").Value = "'0.1845"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").Value = "'10.17"
$ws.Range("E44").Value = "  -8.78%  "
$ws.Range("D45").Value = "'1.253"
$ws.Range("E45").Value = "  -6.10%  "
$ws.Range("D46").Value = "'0.5604"
$ws.Range("E46").Value = "  -6.18%  "
$ws.Range("D47").Value = "'12.18"
$ws.Range("E47").Value = "  -6.41%  "
$ws.Range("D48").Value = "'3.350"
$ws.Range("E48").Value = "  -1.42%  "
$ws.Range("D49").Value = "'1.905"
$ws.Range("E49").Value = "  -7.50%  "
$ws.Range("D50").Value = "'0.06833"
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  -0.66%  "
